$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68 (1-indexed), shifting existing rows 68-97 down to 69-98
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new entry
$ws.Cells.Item(68, 1).Value = "L_FATF_1"
$ws.Cells.Item(68, 2).Value = "Q_FATF"
$ws.Cells.Item(68, 3).Value = "Ratings Deutschland"
$ws.Cells.Item(68, 4).Value = "Ratings Germany"
$ws.Cells.Item(68, 5).Value = "https://www.fatf-gafi.org/en/countries/detail/Germany.html"
$ws.Cells.Item(68, 6).Value = "https://www.fatf-gafi.org/en/countries/detail/Germany.html"
